$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORTDATE: 2020-09-30 -> 2019-09-30
$ws.Range("H2").Value = "2019-09-30 00:00:00"

# BASIC_EPS
$ws.Range("I2").Value = 0.12

# TOTAL_OPERATE_INCOME
$ws.Range("K2").Value = 44576021.69

# PARENT_NETPROFIT
$ws.Range("L2").Value = 3588039.88

# YSTZ, SJLTZ, BPS, MGJYXJJE now blank
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = ""
$ws.Range("Q2").Value = ""

# XSMLL
$ws.Range("R2").Value = 36.4240691619

# ISNEW - numeric-looking text, write via a scratch cell + values-only
# paste so Excel keeps it as text (quote-prefixed) instead of coercing
# to a number, without disturbing AB2's existing (default) cell style.
$ws.Range("ZZ1").Value = "'0"
$ws.Range("ZZ1").Copy()
$ws.Range("AB2").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

# QDATE
$ws.Range("AC2").Value = "2019Q3"

# DATATYPE
$ws.Range("AD2").Value = "2019年 三季报"

# DATAYEAR - numeric-looking text, same scratch-cell technique as AB2
$ws.Range("ZZ1").Value = "'2019"
$ws.Range("ZZ1").Copy()
$ws.Range("AE2").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$excel.CutCopyMode = 0
